# Insert a new data row at row 17 (weekly update), pushing the
# existing rows 17-38 down to 18-39. The new row reuses the same
# constant dimension/category values as the other rows in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Insert()

$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 44467
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 300000001
$ws.Cells.Item(17, 7).Value = "Rabanito"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 20
$ws.Cells.Item(17, 11).Value = 6000
$ws.Cells.Item(17, 12).Value = 6000
$ws.Cells.Item(17, 13).Value = 6000
$ws.Cells.Item(17, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(17, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(17, 16).Value = 500
$ws.Cells.Item(17, 17).Value = 12
$ws.Cells.Item(17, 18).Value = "Hortaliza"
